$d = $word.ActiveDocument

$map = @(
    @("348×6=", "863×3="),
    @("966×2=", "552×8="),
    @("456×9=", "734×5="),
    @("209×8=", "491×5="),
    @("333×6=", "574×6="),
    @("559×9=", "344×4="),
    @("327×3=", "979×7="),
    @("511×8=", "221×7="),
    @("376×9=", "766×2="),
    @("871×9=", "608×4="),
    @("272×6=", "115×5="),
    @("983×4=", "611×6="),
    @("843×5=", "897×6="),
    @("497×7=", "567×3="),
    @("689×6=", "323×5="),
    @("962×8=", "138×9="),
    @("489×3=", "924×6="),
    @("287×8=", "167×5="),
    @("647×6=", "267×4="),
    @("885×8=", "774×2="),
    @("837×4=", "853×4="),
    @("413×9=", "159×4="),
    @("926×9=", "873×6="),
    @("178×6=", "345×6="),
    @("528×4=", "427×4=")
)

foreach ($pair in $map) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
